# Group14-A2 WaiterData.xlsx edit:
# "changed the layout of the page and also changed the active and inactive pages"
#
# - Adds a new "Active Status" column (I1) to the "User Data" sheet header row.
# - Moves the selection/active cell to I5 (reflecting the new layout).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for the "Active Status" column (extends the used range to A1:I1).
$ws.Range("I1").Value = "Active Status"

# Update the sheet's selection to match the new layout (activeCell I5, sqref I5).
$ws.Range("I5").Select() | Out-Null
